$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password for user Gunawan (row 2) from "gungun123" to "gunawan123"
$ws.Range("D2").Value = "gunawan123"

# Move the active selection to D3 (matches the saved sheet view state)
$ws.Range("D3").Select()
